# Update the PSSM matrix (B2:K21, 20 amino-acid rows x 10 position columns)
# with the supplemental-figure recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[,]' 20,10
$values[0,0] = -16.6175105113293
$values[0,1] = 0.5335555915972374
$values[0,2] = -16.6175105113293
$values[0,3] = -16.6175105113293
$values[0,4] = -16.6175105113293
$values[0,5] = -16.6175105113293
$values[0,6] = -16.6175105113293
$values[0,7] = -16.6175105113293
$values[0,8] = -16.6175105113293
$values[0,9] = -16.6175105113293
$values[1,0] = -16.6175105113293
$values[1,1] = -16.6175105113293
$values[1,2] = -16.6175105113293
$values[1,3] = -16.6175105113293
$values[1,4] = -16.6175105113293
$values[1,5] = -16.6175105113293
$values[1,6] = -16.6175105113293
$values[1,7] = 1.21629992582779
$values[1,8] = -16.6175105113293
$values[1,9] = -16.6175105113293
$values[2,0] = -16.6175105113293
$values[2,1] = 0.3847880964802573
$values[2,2] = -16.6175105113293
$values[2,3] = -16.6175105113293
$values[2,4] = 3.876553030975033
$values[2,5] = -16.6175105113293
$values[2,6] = 1.537609116368276
$values[2,7] = -16.6175105113293
$values[2,8] = 2.357097331634492
$values[2,9] = -16.6175105113293
$values[3,0] = -16.6175105113293
$values[3,1] = 0.5200110235821818
$values[3,2] = -16.6175105113293
$values[3,3] = -16.6175105113293
$values[3,4] = -16.6175105113293
$values[3,5] = 3.579424137823911
$values[3,6] = -16.6175105113293
$values[3,7] = -16.6175105113293
$values[3,8] = -16.6175105113293
$values[3,9] = -16.6175105113293
$values[4,0] = -16.6175105113293
$values[4,1] = -16.6175105113293
$values[4,2] = -16.6175105113293
$values[4,3] = -16.6175105113293
$values[4,4] = -16.6175105113293
$values[4,5] = -16.6175105113293
$values[4,6] = -16.6175105113293
$values[4,7] = -16.6175105113293
$values[4,8] = -16.6175105113293
$values[4,9] = -16.6175105113293
$values[5,0] = 2.389168942994096
$values[5,1] = -16.6175105113293
$values[5,2] = -16.6175105113293
$values[5,3] = -16.6175105113293
$values[5,4] = -16.6175105113293
$values[5,5] = -16.6175105113293
$values[5,6] = -16.6175105113293
$values[5,7] = -16.6175105113293
$values[5,8] = -16.6175105113293
$values[5,9] = -16.6175105113293
$values[6,0] = -16.6175105113293
$values[6,1] = -16.6175105113293
$values[6,2] = -16.6175105113293
$values[6,3] = 3.712300057951267
$values[6,4] = -16.6175105113293
$values[6,5] = -16.6175105113293
$values[6,6] = -16.6175105113293
$values[6,7] = -16.6175105113293
$values[6,8] = -16.6175105113293
$values[6,9] = -16.6175105113293
$values[7,0] = 3.883744542181849
$values[7,1] = -16.6175105113293
$values[7,2] = -16.6175105113293
$values[7,3] = -16.6175105113293
$values[7,4] = -16.6175105113293
$values[7,5] = -16.6175105113293
$values[7,6] = -16.6175105113293
$values[7,7] = -16.6175105113293
$values[7,8] = -16.6175105113293
$values[7,9] = -16.6175105113293
$values[8,0] = -16.6175105113293
$values[8,1] = -16.6175105113293
$values[8,2] = -16.6175105113293
$values[8,3] = -16.6175105113293
$values[8,4] = -16.6175105113293
$values[8,5] = -16.6175105113293
$values[8,6] = -16.6175105113293
$values[8,7] = 0.8274153923390818
$values[8,8] = -16.6175105113293
$values[8,9] = 1.76739547015854
$values[9,0] = -16.6175105113293
$values[9,1] = -16.6175105113293
$values[9,2] = -16.6175105113293
$values[9,3] = 1.583048102543625
$values[9,4] = -16.6175105113293
$values[9,5] = 1.492427239496996
$values[9,6] = -16.6175105113293
$values[9,7] = -16.6175105113293
$values[9,8] = -16.6175105113293
$values[9,9] = 1.056881780033542
$values[10,0] = -16.6175105113293
$values[10,1] = -16.6175105113293
$values[10,2] = -16.6175105113293
$values[10,3] = -16.6175105113293
$values[10,4] = -16.6175105113293
$values[10,5] = -16.6175105113293
$values[10,6] = -16.6175105113293
$values[10,7] = -16.6175105113293
$values[10,8] = -16.6175105113293
$values[10,9] = -16.6175105113293
$values[11,0] = -16.6175105113293
$values[11,1] = -16.6175105113293
$values[11,2] = -16.6175105113293
$values[11,3] = 0.6001916488041391
$values[11,4] = -16.6175105113293
$values[11,5] = -16.6175105113293
$values[11,6] = -16.6175105113293
$values[11,7] = -16.6175105113293
$values[11,8] = 1.040094048118556
$values[11,9] = 2.209860106640013
$values[12,0] = -16.6175105113293
$values[12,1] = -16.6175105113293
$values[12,2] = 4.321914463781873
$values[12,3] = -16.6175105113293
$values[12,4] = -16.6175105113293
$values[12,5] = -16.6175105113293
$values[12,6] = -16.6175105113293
$values[12,7] = -16.6175105113293
$values[12,8] = -16.6175105113293
$values[12,9] = 1.677760689227849
$values[13,0] = -16.6175105113293
$values[13,1] = -16.6175105113293
$values[13,2] = -16.6175105113293
$values[13,3] = -16.6175105113293
$values[13,4] = -16.6175105113293
$values[13,5] = -16.6175105113293
$values[13,6] = -16.6175105113293
$values[13,7] = -16.6175105113293
$values[13,8] = -16.6175105113293
$values[13,9] = -16.6175105113293
$values[14,0] = -16.6175105113293
$values[14,1] = -16.6175105113293
$values[14,2] = -16.6175105113293
$values[14,3] = -16.6175105113293
$values[14,4] = -16.6175105113293
$values[14,5] = -16.6175105113293
$values[14,6] = -16.6175105113293
$values[14,7] = -16.6175105113293
$values[14,8] = 2.761522998093269
$values[14,9] = -16.6175105113293
$values[15,0] = -16.6175105113293
$values[15,1] = 0.6043503365586894
$values[15,2] = -16.6175105113293
$values[15,3] = -16.6175105113293
$values[15,4] = -16.6175105113293
$values[15,5] = -16.6175105113293
$values[15,6] = 2.847603465848622
$values[15,7] = 0.195782184811463
$values[15,8] = 1.723056799343581
$values[15,9] = -16.6175105113293
$values[16,0] = -16.6175105113293
$values[16,1] = -16.6175105113293
$values[16,2] = -16.6175105113293
$values[16,3] = -16.6175105113293
$values[16,4] = -16.6175105113293
$values[16,5] = -16.6175105113293
$values[16,6] = 1.862084257727185
$values[16,7] = -0.2676591463279469
$values[16,8] = 1.452930181132768
$values[16,9] = -16.6175105113293
$values[17,0] = -16.6175105113293
$values[17,1] = -16.6175105113293
$values[17,2] = -16.6175105113293
$values[17,3] = -16.6175105113293
$values[17,4] = -16.6175105113293
$values[17,5] = -16.6175105113293
$values[17,6] = 1.417702057938307
$values[17,7] = 2.336165896149824
$values[17,8] = -16.6175105113293
$values[17,9] = -16.6175105113293
$values[18,0] = -16.6175105113293
$values[18,1] = 0.9317402053158966
$values[18,2] = -16.6175105113293
$values[18,3] = -16.6175105113293
$values[18,4] = 2.409242548505295
$values[18,5] = -16.6175105113293
$values[18,6] = 0.8865092744802312
$values[18,7] = 3.14995141908672
$values[18,8] = -16.6175105113293
$values[18,9] = 2.741880886922893
$values[19,0] = -16.6175105113293
$values[19,1] = 3.630506213093189
$values[19,2] = -16.6175105113293
$values[19,3] = 1.251318070735246
$values[19,4] = -16.6175105113293
$values[19,5] = 2.387414909586284
$values[19,6] = 0.8016349588265542
$values[19,7] = -16.6175105113293
$values[19,8] = -16.6175105113293
$values[19,9] = -16.6175105113293

$ws.Range("B2:K21").Value = $values
